$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45875.01041666666
$ws.Cells.Item(2, 2).Value = 65
$ws.Cells.Item(3, 1).Value = 45875.02083333334
$ws.Cells.Item(3, 2).Value = 65
$ws.Cells.Item(4, 1).Value = 45875.03125
$ws.Cells.Item(4, 2).Value = 65
$ws.Cells.Item(5, 1).Value = 45875.04166666666
$ws.Cells.Item(5, 2).Value = 65
$ws.Cells.Item(6, 1).Value = 45875.05208333334
$ws.Cells.Item(6, 2).Value = 56
$ws.Cells.Item(7, 1).Value = 45875.0625
$ws.Cells.Item(7, 2).Value = 56
$ws.Cells.Item(8, 1).Value = 45875.07291666666
$ws.Cells.Item(8, 2).Value = 57
$ws.Cells.Item(9, 1).Value = 45875.08333333334
$ws.Cells.Item(9, 2).Value = 57
$ws.Cells.Item(10, 1).Value = 45875.09375
$ws.Cells.Item(10, 2).Value = 57
$ws.Cells.Item(11, 1).Value = 45875.10416666666
$ws.Cells.Item(11, 2).Value = 58
$ws.Cells.Item(12, 1).Value = 45875.11458333334
$ws.Cells.Item(12, 2).Value = 58
$ws.Cells.Item(13, 1).Value = 45875.125
$ws.Cells.Item(13, 2).Value = 59
$ws.Cells.Item(14, 1).Value = 45875.13541666666
$ws.Cells.Item(14, 2).Value = 69
$ws.Cells.Item(15, 1).Value = 45875.14583333334
$ws.Cells.Item(15, 2).Value = 69
$ws.Cells.Item(16, 1).Value = 45875.15625
$ws.Cells.Item(16, 2).Value = 69
$ws.Cells.Item(17, 1).Value = 45875.16666666666
$ws.Cells.Item(17, 2).Value = 69
$ws.Cells.Item(18, 1).Value = 45875.17708333334
$ws.Cells.Item(18, 2).Value = 83
$ws.Cells.Item(19, 1).Value = 45875.1875
$ws.Cells.Item(19, 2).Value = 83
$ws.Cells.Item(20, 1).Value = 45875.19791666666
$ws.Cells.Item(20, 2).Value = 82
$ws.Cells.Item(21, 1).Value = 45875.20833333334
$ws.Cells.Item(21, 2).Value = 82
$ws.Cells.Item(22, 1).Value = 45875.21875
$ws.Cells.Item(22, 2).Value = 89
$ws.Cells.Item(23, 1).Value = 45875.22916666666
$ws.Cells.Item(23, 2).Value = 89
$ws.Cells.Item(24, 1).Value = 45875.23958333334
$ws.Cells.Item(24, 2).Value = 89
$ws.Cells.Item(25, 1).Value = 45875.25
$ws.Cells.Item(25, 2).Value = 89
$ws.Cells.Item(26, 1).Value = 45875.26041666666
$ws.Cells.Item(26, 2).Value = 82
$ws.Cells.Item(27, 1).Value = 45875.27083333334
$ws.Cells.Item(27, 2).Value = 82
$ws.Cells.Item(28, 1).Value = 45875.28125
$ws.Cells.Item(28, 2).Value = 82
$ws.Cells.Item(29, 1).Value = 45875.29166666666
$ws.Cells.Item(29, 2).Value = 83
$ws.Cells.Item(30, 1).Value = 45875.30208333334
$ws.Cells.Item(30, 2).Value = 66
$ws.Cells.Item(31, 1).Value = 45875.3125
$ws.Cells.Item(31, 2).Value = 67
$ws.Cells.Item(32, 1).Value = 45875.32291666666
$ws.Cells.Item(32, 2).Value = 67
$ws.Cells.Item(33, 1).Value = 45875.33333333334
$ws.Cells.Item(33, 2).Value = 67
$ws.Cells.Item(34, 1).Value = 45875.34375
$ws.Cells.Item(34, 2).Value = 57
$ws.Cells.Item(35, 1).Value = 45875.35416666666
$ws.Cells.Item(35, 2).Value = 57
$ws.Cells.Item(36, 1).Value = 45875.36458333334
$ws.Cells.Item(36, 2).Value = 58
$ws.Cells.Item(37, 1).Value = 45875.375
$ws.Cells.Item(37, 2).Value = 58
$ws.Cells.Item(38, 1).Value = 45875.38541666666
$ws.Cells.Item(38, 2).Value = 63
$ws.Cells.Item(39, 1).Value = 45875.39583333334
$ws.Cells.Item(39, 2).Value = 63
$ws.Cells.Item(40, 1).Value = 45875.40625
$ws.Cells.Item(40, 2).Value = 63
$ws.Cells.Item(41, 1).Value = 45875.41666666666
$ws.Cells.Item(41, 2).Value = 64
$ws.Cells.Item(42, 1).Value = 45875.42708333334
$ws.Cells.Item(42, 2).Value = 77
$ws.Cells.Item(43, 1).Value = 45875.4375
$ws.Cells.Item(43, 2).Value = 78
$ws.Cells.Item(44, 1).Value = 45875.44791666666
$ws.Cells.Item(44, 2).Value = 78
$ws.Cells.Item(45, 1).Value = 45875.45833333334
$ws.Cells.Item(45, 2).Value = 78
$ws.Cells.Item(46, 1).Value = 45875.46875
$ws.Cells.Item(46, 2).Value = 100
$ws.Cells.Item(47, 1).Value = 45875.47916666666
$ws.Cells.Item(47, 2).Value = 100
$ws.Cells.Item(48, 1).Value = 45875.48958333334
$ws.Cells.Item(48, 2).Value = 101
$ws.Cells.Item(49, 1).Value = 45875.5
$ws.Cells.Item(49, 2).Value = 102
$ws.Cells.Item(50, 1).Value = 45875.51041666666
$ws.Cells.Item(50, 2).Value = 125
$ws.Cells.Item(51, 1).Value = 45875.52083333334
$ws.Cells.Item(51, 2).Value = 125
$ws.Cells.Item(52, 1).Value = 45875.53125
$ws.Cells.Item(52, 2).Value = 126
$ws.Cells.Item(53, 1).Value = 45875.54166666666
$ws.Cells.Item(53, 2).Value = 127
$ws.Cells.Item(54, 1).Value = 45875.55208333334
$ws.Cells.Item(54, 2).Value = 171
$ws.Cells.Item(55, 1).Value = 45875.5625
$ws.Cells.Item(55, 2).Value = 172
$ws.Cells.Item(56, 1).Value = 45875.57291666666
$ws.Cells.Item(56, 2).Value = 173
$ws.Cells.Item(57, 1).Value = 45875.58333333334
$ws.Cells.Item(57, 2).Value = 174
$ws.Cells.Item(58, 1).Value = 45875.59375
$ws.Cells.Item(58, 2).Value = 238
$ws.Cells.Item(59, 1).Value = 45875.60416666666
$ws.Cells.Item(59, 2).Value = 239
$ws.Cells.Item(60, 1).Value = 45875.61458333334
$ws.Cells.Item(60, 2).Value = 239
$ws.Cells.Item(61, 1).Value = 45875.625
$ws.Cells.Item(61, 2).Value = 240
$ws.Cells.Item(62, 1).Value = 45875.63541666666
$ws.Cells.Item(62, 2).Value = 312
$ws.Cells.Item(63, 1).Value = 45875.64583333334
$ws.Cells.Item(63, 2).Value = 312
$ws.Cells.Item(64, 1).Value = 45875.65625
$ws.Cells.Item(64, 2).Value = 312
$ws.Cells.Item(65, 1).Value = 45875.66666666666
$ws.Cells.Item(65, 2).Value = 311
$ws.Cells.Item(66, 1).Value = 45875.67708333334
$ws.Cells.Item(66, 2).Value = 350
$ws.Cells.Item(67, 1).Value = 45875.6875
$ws.Cells.Item(67, 2).Value = 350
$ws.Cells.Item(68, 1).Value = 45875.69791666666
$ws.Cells.Item(68, 2).Value = 351
$ws.Cells.Item(69, 1).Value = 45875.70833333334
$ws.Cells.Item(69, 2).Value = 351
$ws.Cells.Item(70, 1).Value = 45875.71875
$ws.Cells.Item(70, 2).Value = 368
$ws.Cells.Item(71, 1).Value = 45875.72916666666
$ws.Cells.Item(71, 2).Value = 369
$ws.Cells.Item(72, 1).Value = 45875.73958333334
$ws.Cells.Item(72, 2).Value = 370
$ws.Cells.Item(73, 1).Value = 45875.75
$ws.Cells.Item(73, 2).Value = 371
$ws.Cells.Item(74, 1).Value = 45875.76041666666
$ws.Cells.Item(74, 2).Value = 371
$ws.Cells.Item(75, 1).Value = 45875.77083333334
$ws.Cells.Item(75, 2).Value = 372
$ws.Cells.Item(76, 1).Value = 45875.78125
$ws.Cells.Item(76, 2).Value = 374
$ws.Cells.Item(77, 1).Value = 45875.79166666666
$ws.Cells.Item(77, 2).Value = 376
$ws.Cells.Item(78, 1).Value = 45875.80208333334
$ws.Cells.Item(78, 2).Value = 391
$ws.Cells.Item(79, 1).Value = 45875.8125
$ws.Cells.Item(79, 2).Value = 394
$ws.Cells.Item(80, 1).Value = 45875.82291666666
$ws.Cells.Item(80, 2).Value = 398
$ws.Cells.Item(81, 1).Value = 45875.83333333334
$ws.Cells.Item(81, 2).Value = 401
$ws.Cells.Item(82, 1).Value = 45875.84375
$ws.Cells.Item(82, 2).Value = 432
$ws.Cells.Item(83, 1).Value = 45875.85416666666
$ws.Cells.Item(83, 2).Value = 437
$ws.Cells.Item(84, 1).Value = 45875.86458333334
$ws.Cells.Item(84, 2).Value = 442
$ws.Cells.Item(85, 1).Value = 45875.875
$ws.Cells.Item(85, 2).Value = 447
$ws.Cells.Item(86, 1).Value = 45875.88541666666
$ws.Cells.Item(86, 2).Value = 618
$ws.Cells.Item(87, 1).Value = 45875.89583333334
$ws.Cells.Item(87, 2).Value = 621
$ws.Cells.Item(88, 1).Value = 45875.90625
$ws.Cells.Item(88, 2).Value = 623
$ws.Cells.Item(89, 1).Value = 45875.91666666666
$ws.Cells.Item(89, 2).Value = 626
$ws.Cells.Item(90, 1).Value = 45875.92708333334
$ws.Cells.Item(90, 2).Value = 852
$ws.Cells.Item(91, 1).Value = 45875.9375
$ws.Cells.Item(91, 2).Value = 857
$ws.Cells.Item(92, 1).Value = 45875.94791666666
$ws.Cells.Item(92, 2).Value = 862
$ws.Cells.Item(93, 1).Value = 45875.95833333334
$ws.Cells.Item(93, 2).Value = 866
$ws.Cells.Item(94, 1).Value = 45875.96875
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(95, 1).Value = 45875.97916666666
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(96, 1).Value = 45875.98958333334
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(97, 1).Value = 45876
$ws.Cells.Item(97, 2).Value = 0
